$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.679.26"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "1.961.89"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.58"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.30"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +7.58%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.376"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0796"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -6.12%  "
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("E12").Value = "  +6.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.00"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.834"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.08%  "
$ws.Range("D15").Value = "2.251.20"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.29"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.84%  "
$ws.Range("D17").Value = "1.963.67"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "36.596.73"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.87"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "230.22"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  +4.37%  "
$ws.Range("E25").Value = "  +2.77%  "
$ws.Range("E26").Value = "  +8.52%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.90"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.41"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.31"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +19.59%  "
$ws.Range("E31").Value = "  +1.64%  "
$ws.Range("E32").Value = "  +4.94%  "
$ws.Range("E33").Value = "  +0.21%  "
$ws.Range("E34").Value = "  +6.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.60"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +13.82%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.28"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.52%  "
$ws.Range("B37").Value = "BinanceUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.16%  "
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("E39").Value = "  -10.92%  "
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.17"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.05"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("D45").Value = "1.370.31"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.87"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.60%  "
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.17"
$ws.Range("D48").ClearFormats()
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.41"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.37%  "
$ws.Range("D51").Value = "2.140.57"
$ws.Range("E51").Value = "  +1.11%  "
